$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1, reusing the same formatting as the neighboring
# header cell (G1) so it shares the existing style index.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the "Save" column values for rows 2-5 (plain, unstyled numbers)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
